# Apply cryptos.xlsx price/volume updates per commit diff (2024-05-19)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.696.12"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "3.073.47"
$ws.Range("E3").Value = "  -1.23%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.24"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.84"
$ws.Range("E6").Value = "  -1.24%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.069.44"
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.43"
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("E11").Value = "  -1.58%  "
$ws.Range("E12").Value = "  -3.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000239"
$ws.Range("E13").Value = "  -2.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.65"
$ws.Range("E14").Value = "  -3.67%  "
$ws.Range("E15").Value = "  -1.73%  "
$ws.Range("D16").Value = "3.587.32"
$ws.Range("E16").Value = "  -1.13%  "
$ws.Range("D17").Value = "66.695.77"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.01"
$ws.Range("E18").Value = "  +4.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.97"
$ws.Range("E19").Value = "  -2.83%  "
$ws.Range("D20").Value = "3.064.61"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "487.75"
$ws.Range("E21").Value = "  +2.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.69"
$ws.Range("E22").Value = "  -2.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.686"
$ws.Range("E23").Value = "  -3.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.64"
$ws.Range("E24").Value = "  -1.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.67"
$ws.Range("E25").Value = "  -5.10%  "
$ws.Range("E26").Value = "  -3.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.17"
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.77"
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.26"
$ws.Range("E30").Value = "  -4.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.57"
$ws.Range("E32").Value = "  -3.45%  "
$ws.Range("E33").Value = "  -3.87%  "
$ws.Range("D34").Value = "0.0₃0910"
$ws.Range("E34").Value = "  -3.25%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  -4.38%  "
$ws.Range("E37").Value = "  -2.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.00"
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("E40").Value = "  -4.67%  "
$ws.Range("E41").Value = "  -3.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.27"
$ws.Range("E42").Value = "  -5.02%  "
$ws.Range("D43").Value = "2.755.90"
$ws.Range("E43").Value = "  -2.64%  "
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "135.68"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0344"
$ws.Range("E46").Value = "  -3.14%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "367.65"
$ws.Range("E47").Value = "  -4.85%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.71"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.15"
$ws.Range("E50").Value = "  -1.99%  "
$ws.Range("E51").Value = "  -1.95%  "
